$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "version 2"
$ws.Range("B2").Value = "cams taba ganda"
